{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the text changes described by the diff: updates the title,\n// rewrites the \"What we like\" / \"What we don't like\" bullet lists, and\n// updates the closing bold title + italic meta-description paragraph.\n\nconst replacements = [\n  {\n    find: \"Play Hotline 2 Free - Review of Gameplay, Design, and Bonus Features\",\n    replace: \"Play Hotline 2 Free - Exciting and Challenging Slot Game\",\n  },\n  {\n    find: \"Unique hotlines feature\",\n    replace: \"Traditional slot game structure with 5 reels and 243 ways to win\",\n  },\n  {\n    find: \"Design and graphics are classic and modern\",\n    replace: \"Unique hotline feature for increased chances of winning\",\n  },\n  {\n    find: \"Moderately volatile with potential for big wins\",\n    replace: \"Classic yet modern design with attractive graphics\",\n  },\n  {\n    find: \"Free spins and Wild/Scatter symbols for bonus features\",\n    replace: \"Moderate volatility providing consistent small wins\",\n  },\n  {\n    find: \"Patience required for substantial payouts\",\n    replace: \"Substantial payouts may require more patience\",\n  },\n  {\n    find: \"Typical RTP rate for slot games\",\n    replace: \"Limited number of free spins\",\n  },\n  {\n    find: \"Read our review of Hotline 2, a classic/modern slot game with free spins, Wild/Scatter symbols, and a unique hotlines feature. Play it for free now.\",\n    replace: \"Read our Hotline 2 review and play this exciting and challenging slot game for free.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the text changes described by the diff: updates the title,\n# rewrites the \"What we like\" / \"What we don't like\" bullet lists, and\n# updates the closing bold title + italic meta-description paragraph.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @{ Find = \"Play Hotline 2 Free - Review of Gameplay, Design, and Bonus Features\"; Replace = \"Play Hotline 2 Free - Exciting and Challenging Slot Game\" },\n    @{ Find = \"Unique hotlines feature\"; Replace = \"Traditional slot game structure with 5 reels and 243 ways to win\" },\n    @{ Find = \"Design and graphics are classic and modern\"; Replace = \"Unique hotline feature for increased chances of winning\" },\n    @{ Find = \"Moderately volatile with potential for big wins\"; Replace = \"Classic yet modern design with attractive graphics\" },\n    @{ Find = \"Free spins and Wild/Scatter symbols for bonus features\"; Replace = \"Moderate volatility providing consistent small wins\" },\n    @{ Find = \"Patience required for substantial payouts\"; Replace = \"Substantial payouts may require more patience\" },\n    @{ Find = \"Typical RTP rate for slot games\"; Replace = \"Limited number of free spins\" },\n    @{ Find = \"Read our review of Hotline 2, a classic/modern slot game with free spins, Wild/Scatter symbols, and a unique hotlines feature. Play it for free now.\"; Replace = \"Read our Hotline 2 review and play this exciting and challenging slot game for free.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n}\n"}
